$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 5199.75
$ws.Range("I74").Value = 5266.3335
$ws.Range("K74").Value = 5266.3335
$ws.Range("M74").Value = -4330.3335
$ws.Range("H77").Value = 5199.75
$ws.Range("I77").Value = 5266.3335
$ws.Range("K77").Value = 26331.6675
$ws.Range("M77").Value = -21651.6675
$ws.Range("H131").Value = 4849.8
$ws.Range("I131").Value = 5312.25
$ws.Range("J131").Value = 3000
$ws.Range("K131").Value = 15936.75
$ws.Range("L131").Value = 9000
$ws.Range("M131").Value = -10896.75
$ws.Range("N131").Value = -19080
$ws.Range("H135").Value = 2260.2666
$ws.Range("I135").Value = 1993.5116
$ws.Range("K135").Value = 17941.6044
$ws.Range("M135").Value = -15406.6044
$ws.Range("H137").Value = 4856.857
$ws.Range("I137").Value = 4856.857
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 14570.571
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -12020.571
$ws.Range("N137").ClearContents()
$ws.Range("H138").Value = 2970.785
$ws.Range("J138").Value = 3972.9333
$ws.Range("L138").Value = 11918.7999
$ws.Range("N138").Value = -22198.7999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4610.0435
$ws.Range("I2").Value = 663.6316
$ws.Range("J2").Value = 23355.5
$ws.Range("K2").Value = 663.6316
$ws.Range("L2").Value = 23355.5
$ws.Range("M2").Value = -550.6316
$ws.Range("N2").Value = -23581.5
$ws.Range("H32").Value = 2444.5361
$ws.Range("I32").Value = 2114.5522
$ws.Range("K32").Value = 2114.5522
$ws.Range("M32").Value = -1827.5522
$ws.Range("H74").Value = 3678.0588
$ws.Range("J74").Value = 4842.778
$ws.Range("L74").Value = 4842.778
$ws.Range("N74").Value = -6590.778
$ws.Range("H77").Value = 3678.0588
$ws.Range("J77").Value = 4842.778
$ws.Range("L77").Value = 24213.89
$ws.Range("N77").Value = -32949.89
$ws.Range("H116").Value = 4610.0435
$ws.Range("I116").Value = 663.6316
$ws.Range("J116").Value = 23355.5
$ws.Range("K116").Value = 663.6316
$ws.Range("L116").Value = 23355.5
$ws.Range("M116").Value = 1630.3684
$ws.Range("N116").Value = -27943.5
$ws.Range("H122").Value = 4260.533
$ws.Range("I122").Value = 4411.154
$ws.Range("K122").Value = 13233.462
$ws.Range("M122").Value = -10783.462
$ws.Range("H132").Value = 9039.833000000001
$ws.Range("I132").Value = 5586.6665
$ws.Range("K132").Value = 16759.9995
$ws.Range("M132").Value = -14229.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4610.0435
$ws.Range("I3").Value = 663.6316
$ws.Range("J3").Value = 23355.5
$ws.Range("K3").Value = 663.6316
$ws.Range("L3").Value = 23355.5
$ws.Range("M3").Value = -549.6316
$ws.Range("N3").Value = -23583.5
$ws.Range("H20").Value = 52633416
$ws.Range("I20").Value = 90911120
$ws.Range("K20").Value = 90911120
$ws.Range("M20").Value = -90910873
$ws.Range("H80").Value = 801.1818
$ws.Range("I80").Value = 717.1429000000001
$ws.Range("J80").Value = 840.4
$ws.Range("K80").Value = 717.1429000000001
$ws.Range("L80").Value = 840.4
$ws.Range("M80").Value = 280.8570999999999
$ws.Range("N80").Value = -2836.4
$ws.Range("H83").Value = 801.1818
$ws.Range("I83").Value = 717.1429000000001
$ws.Range("J83").Value = 840.4
$ws.Range("K83").Value = 3585.7145
$ws.Range("L83").Value = 4202
$ws.Range("M83").Value = 1406.2855
$ws.Range("N83").Value = -14186
$ws.Range("H134").Value = 3288.7727
$ws.Range("I134").Value = 2005.25
$ws.Range("K134").Value = 6015.75
$ws.Range("M134").Value = -3480.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6807.3335
$ws.Range("I31").Value = 5211.5
$ws.Range("J31").Value = 9999
$ws.Range("K31").Value = 5211.5
$ws.Range("L31").Value = 9999
$ws.Range("M31").Value = -4916.5
$ws.Range("N31").Value = -10589
$ws.Range("H34").Value = 6807.3335
$ws.Range("I34").Value = 5211.5
$ws.Range("J34").Value = 9999
$ws.Range("K34").Value = 5211.5
$ws.Range("L34").Value = 9999
$ws.Range("M34").Value = -5009.5
$ws.Range("N34").Value = -10403
$ws.Range("H58").Value = 3340.5
$ws.Range("I58").Value = 3534.318
$ws.Range("J58").Value = 2274.5
$ws.Range("K58").Value = 3534.318
$ws.Range("L58").Value = 2274.5
$ws.Range("M58").Value = -3331.318
$ws.Range("N58").Value = -2680.5
$ws.Range("H132").Value = 1154.575
$ws.Range("I132").Value = 1093.4517
$ws.Range("J132").Value = 1365.1111
$ws.Range("K132").Value = 3280.3551
$ws.Range("L132").Value = 4095.3333
$ws.Range("M132").Value = -750.3551000000002
$ws.Range("N132").Value = -9155.3333
$ws.Range("H134").Value = 1724.7435
$ws.Range("I134").Value = 1719.0286
$ws.Range("J134").Value = 1774.75
$ws.Range("K134").Value = 5157.085800000001
$ws.Range("L134").Value = 5324.25
$ws.Range("M134").Value = -2622.085800000001
$ws.Range("N134").Value = -10394.25
$ws.Range("H136").Value = 3340.5
$ws.Range("I136").Value = 3534.318
$ws.Range("J136").Value = 2274.5
$ws.Range("K136").Value = 10602.954
$ws.Range("L136").Value = 6823.5
$ws.Range("M136").Value = -8052.954000000002
$ws.Range("N136").Value = -11923.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 996.4286
$ws.Range("I47").Value = 495.83334
$ws.Range("K47").Value = 1487.50002
$ws.Range("M47").Value = -1056.50002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1674.2222
$ws.Range("I113").Value = 1361.25
$ws.Range("J113").Value = 1924.6
$ws.Range("K113").Value = 1361.25
$ws.Range("L113").Value = 1924.6
$ws.Range("M113").Value = 808.75
$ws.Range("N113").Value = -6264.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3968.7896
$ws.Range("I7").Value = 3259.2354
$ws.Range("J7").Value = 10000
$ws.Range("K7").Value = 3259.2354
$ws.Range("L7").Value = 10000
$ws.Range("M7").Value = -3147.2354
$ws.Range("N7").Value = -10224
$ws.Range("H46").Value = 808.1667
$ws.Range("I46").Value = 899.6667
$ws.Range("J46").Value = 716.6667
$ws.Range("K46").Value = 899.6667
$ws.Range("L46").Value = 716.6667
$ws.Range("M46").Value = -711.6667
$ws.Range("N46").Value = -1092.6667
$ws.Range("H126").Value = 3968.7896
$ws.Range("I126").Value = 3259.2354
$ws.Range("J126").Value = 10000
$ws.Range("K126").Value = 9777.706200000001
$ws.Range("L126").Value = 30000
$ws.Range("M126").Value = -7307.706200000001
$ws.Range("N126").Value = -34940
$ws.Range("H132").Value = 4785.8857
$ws.Range("I132").Value = 3176.6956
$ws.Range("K132").Value = 9530.086800000001
$ws.Range("M132").Value = -7000.086800000001
$ws.Range("H136").Value = 7191.647
$ws.Range("I136").Value = 6758.5386
$ws.Range("J136").Value = 8599.25
$ws.Range("K136").Value = 20275.6158
$ws.Range("L136").Value = 25797.75
$ws.Range("M136").Value = -17725.6158
$ws.Range("N136").Value = -30897.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2157.6667
$ws.Range("I132").Value = 999.75
$ws.Range("J132").Value = 2736.625
$ws.Range("K132").Value = 2999.25
$ws.Range("L132").Value = 8209.875
$ws.Range("M132").Value = -469.25
$ws.Range("N132").Value = -13269.875
$ws.Range("H136").Value = 5143.718
$ws.Range("I136").Value = 6075.7417
$ws.Range("J136").Value = 1532.125
$ws.Range("K136").Value = 18227.2251
$ws.Range("L136").Value = 4596.375
$ws.Range("M136").Value = -15677.2251
$ws.Range("N136").Value = -9696.375
